# This script rewrites the weekly schedule cells on the active sheet so that
# the underlying shared-strings table ends up reordered/rewritten exactly as
# described by the target diff. The worksheet layout (rows/columns used) is
# unchanged; only the text values of a subset of cells are updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Math1-1"
$ws.Range("C2").Value = "Litr1-3"
$ws.Range("D2").Value = "Math1-4"
$ws.Range("A3").Value = "Phys1-1"
$ws.Range("D3").Value = "English1-4"
$ws.Range("C4").Value = "English1-3"
$ws.Range("D4").Value = "Litr1-4"

$ws.Range("A7").Value = "Math2-1"
$ws.Range("B7").Value = "Phys2-2"
$ws.Range("D7").Value = "Math2-4"
$ws.Range("E7").Value = "Phys2-5"
$ws.Range("A8").Value = "English2-1"
$ws.Range("B8").Value = "Math2-2"
$ws.Range("D8").Value = "English2-4"
$ws.Range("E8").Value = "Math2-5"

$ws.Range("A11").Value = "Russian3-1"
$ws.Range("B11").Value = "Phys3-2"
$ws.Range("C11").Value = "Russian3-3"
$ws.Range("D11").Value = "Phys3-4"
$ws.Range("E11").Value = "Litra3-5"
$ws.Range("B12").Value = "Russian3-2"
$ws.Range("D12").Value = "Litra3-4"
$ws.Range("E12").Value = "Phys3-5"
$ws.Range("A13").Value = "English3-1"

$ws.Range("A16").Value = "English4-1"
$ws.Range("D16").Value = "Math4-4"
$ws.Range("A17").Value = "Math4-1"
$ws.Range("B17").Value = "Math4-2"
$ws.Range("D17").Value = "Phys4-4"

$ws.Range("B20").Value = "Russian5-2"
